$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'41.174.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.88%  "

# Row 3
$ws.Range("D3").Value = "'2.247.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.71%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'302.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.94%  "

# Row 6
$ws.Range("D6").Value = "'91.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.29%  "

# Row 9
$ws.Range("E9").Value = "  +2.22%  "

# Row 10
$ws.Range("D10").Value = "'53.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.22%  "

# Row 11
$ws.Range("D11").Value = "'31.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.80%  "

# Row 12
$ws.Range("D12").Value = "'0.0793"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.07%  "

# Row 13
$ws.Range("E13").Value = "  +3.14%  "

# Row 14
$ws.Range("E14").Value = "  +2.06%  "

# Row 15
$ws.Range("D15").Value = "'2.596.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.79%  "

# Row 16
$ws.Range("E16").Value = "  +2.83%  "

# Row 17
$ws.Range("D17").Value = "'2.244.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.17%  "

# Row 18
$ws.Range("D18").Value = "'0.747"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.18%  "

# Row 19
$ws.Range("D19").Value = "'41.110.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.94%  "

# Row 20
$ws.Range("D20").Value = "'12.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.03%  "

# Row 21
$ws.Range("D21").Value = "'0.0₃0901"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.00%  "

# Row 22
$ws.Range("E22").Value = "  +1.73%  "

# Row 23
$ws.Range("D23").Value = "'66.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.13%  "

# Row 24
$ws.Range("D24").Value = "'240.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.28%  "

# Row 25
$ws.Range("E25").Value = "  +4.36%  "

# Row 26
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("E27").Value = "  +3.18%  "

# Row 28
$ws.Range("D28").Value = "'23.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.56%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.43%  "

# Row 30
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.85%  "

# Row 31
$ws.Range("D31").Value = "'158.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.02%  "

# Row 32
$ws.Range("D32").Value = "'33.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.49%  "

# Row 33
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("D34").Value = "'5.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.43%  "

# Row 35
$ws.Range("D35").Value = "'0.0732"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.34%  "

# Row 36
$ws.Range("E36").Value = "  +7.62%  "

# Row 37
$ws.Range("E37").Value = "  +1.09%  "

# Row 38
$ws.Range("D38").Value = "'16.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.77%  "

# Row 39
$ws.Range("E39").Value = "  +2.79%  "

# Row 40
$ws.Range("E40").Value = "  +5.29%  "

# Row 41
$ws.Range("E41").Value = "  +6.90%  "

# Row 42
$ws.Range("D42").Value = "'3.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.00%  "

# Row 43
$ws.Range("D43").Value = "'20.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.24%  "

# Row 44
$ws.Range("D44").Value = "'2.064.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.75%  "

# Row 45
$ws.Range("D45").Value = "'0.0275"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.26%  "

# Row 46
$ws.Range("D46").Value = "'10.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.15%  "

# Row 47
$ws.Range("E47").Value = "  +11.48%  "

# Row 48
$ws.Range("D48").Value = "'2.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.38%  "

# Row 49
$ws.Range("D49").Value = "'2.468.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.99%  "

# Row 50
$ws.Range("E50").Value = "  +4.01%  "

# Row 51
$ws.Range("E51").Value = "  +1.58%  "
